# The two data rows (Id 60493331 / Id 87015259) simply traded places in the
# sheet: everything that was in row 4 is now in row 5, and vice versa.
# We swap the two rows column by column (cell by cell, to avoid precision
# loss / unwanted type coercion that can happen with bulk array Value2
# assignment), taking care that date-looking text (columns Y and AA) is
# written back as plain text rather than being auto-converted to a real date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 51  # column AY

# Columns whose text must not be auto-interpreted as a date by Excel.
$dateTextCols = @(25, 27)  # Y, AA

for ($col = 1; $col -le $lastCol; $col++) {

    $cell4 = $ws.Cells.Item(4, $col)
    $cell5 = $ws.Cells.Item(5, $col)

    $val4 = $cell4.Value2
    $val5 = $cell5.Value2

    $isDateTextCol = $dateTextCols -contains $col

    if ($isDateTextCol) {
        $cell4.NumberFormat = "@"
        $cell5.NumberFormat = "@"
    }

    if ([string]::IsNullOrEmpty($val5)) {
        $cell4.ClearContents()
    } else {
        $cell4.Value2 = $val5
    }

    if ([string]::IsNullOrEmpty($val4)) {
        $cell5.ClearContents()
    } else {
        $cell5.Value2 = $val4
    }

    if ($isDateTextCol) {
        $cell4.Style = "Normal"
        $cell5.Style = "Normal"
    }
}
